$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column A (data starts at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# 1) Update the "Förändrad" (changed) date in column C for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# 2) Add a friendly display-text second argument to the HYPERLINK formulas
#    in columns S, T, V, W, X, Y (only present for the first few rows).
$linkCols = @("S", "T", "V", "W", "X", "Y")
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Text
    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $r)
        $f = $cell.Formula
        if ($f -ne "" -and $f.IndexOf("HYPERLINK(") -ge 0 -and $f.IndexOf(",") -lt 0) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $name + '")'
            $cell.Formula = $newFormula
        }
    }
}
